# Update the Thpo-Mpl ligand-receptor pair table with new TPM-derived NATMI
# output: refreshed specificity/weight columns and three additional
# Sending-cluster x Target-cluster combinations (rows 5-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Thpo", "Mpl", "ECs", 2, 0.6666666666666666, 0.7896443333333334, 2.368933, 0.6394500843671492, 0.6394500843671492, 3, 1, 0.3210106666666667, 0.9630320000000001, 0.3467150827570242, 0.3467150827570242, 0.253484253872889, 2.281358284856001, 0.2217069889203422, 0.2217069889203423),
    @("FAPs", "Thpo", "Mpl", "FAPs", 2, 0.6666666666666666, 0.7896443333333334, 2.368933, 0.6394500843671492, 0.6394500843671492, 1, 0.3333333333333333, 0.3404013333333333, 1.021204, 0.3676584260666355, 0.3676584260666355, 0.2687959839257778, 2.419163855332, 0.2350992115666033, 0.2350992115666034),
    @("FAPs", "Thpo", "Mpl", "Resolving-Mac", 2, 0.6666666666666666, 0.7896443333333334, 2.368933, 0.6394500843671492, 0.6394500843671492, 1, 0.3333333333333333, 0.264451, 0.793353, 0.2856264911763403, 0.2856264911763404, 0.2088222335943333, 1.879400102349, 0.1826438838802036, 0.1826438838802037),
    @("MuSCs", "Thpo", "Mpl", "ECs", 2, 0.6666666666666666, 0.445236, 1.335708, 0.3605499156328508, 0.3605499156328508, 3, 1, 0.3210106666666667, 0.9630320000000001, 0.3467150827570242, 0.3467150827570242, 0.142925505184, 1.286329546656, 0.1250080938366819, 0.125008093836682),
    @("MuSCs", "Thpo", "Mpl", "FAPs", 2, 0.6666666666666666, 0.445236, 1.335708, 0.3605499156328508, 0.3605499156328508, 1, 0.3333333333333333, 0.3404013333333333, 1.021204, 0.3676584260666355, 0.3676584260666355, 0.151558928048, 1.364030352432, 0.1325592145000321, 0.1325592145000321),
    @("MuSCs", "Thpo", "Mpl", "Resolving-Mac", 2, 0.6666666666666666, 0.445236, 1.335708, 0.3605499156328508, 0.3605499156328508, 1, 0.3333333333333333, 0.264451, 0.793353, 0.2856264911763403, 0.2856264911763404, 0.117743105436, 1.059687948924, 0.1029826072961367, 0.1029826072961367)
)

$startRow = 2
$startCol = 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
